$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 201
$ws1.Range("F6").Value = 510
$ws1.Range("F8").Value = 121
$ws1.Range("F10").Value = 6807
$ws1.Range("F12").Value = 379
$ws1.Range("F13").Value = 3128
$ws1.Range("F15").Value = 364
$ws1.Range("F17").Value = 552
$ws1.Range("F18").Value = 9

# Sheet "演出" (sheet2): update column F values
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 48

# Sheet "全部类型" (sheet4): update column F values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 201
$ws4.Range("F8").Value = 510
$ws4.Range("F10").Value = 121
$ws4.Range("F13").Value = 6807
$ws4.Range("F14").Value = 48
$ws4.Range("F16").Value = 379
$ws4.Range("F17").Value = 3128
$ws4.Range("F19").Value = 364
$ws4.Range("F21").Value = 552
$ws4.Range("F22").Value = 9
